# Terrain/event option table: add the "pick up" (拾取) entry as row 7,
# matching the Id / Name / EventName pattern of the existing rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 1009
$ws.Range("B7").Value = "拾取"
$ws.Range("C7").Value = "OnPickUpItemEvent"

# Copy the formatting from row 6 (same column layout/fonts) onto the new row
# so the new cells share the existing column styles instead of defaults.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)

$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)

$ws.Range("C6").Copy()
$ws.Range("C7").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("C7").Select()
